# Update the "想去人数" (want-to-go count) figures in column F of the
# "展览" and "全部类型" worksheets to the refreshed values captured in the
# latest data scrape.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1235
$ws1.Range("F4").Value  = 23
$ws1.Range("F5").Value  = 12538
$ws1.Range("F7").Value  = 28
$ws1.Range("F8").Value  = 28
$ws1.Range("F10").Value = 12417
$ws1.Range("F11").Value = 240
$ws1.Range("F12").Value = 4896
$ws1.Range("F13").Value = 4823
$ws1.Range("F14").Value = 153
$ws1.Range("F16").Value = 426
$ws1.Range("F17").Value = 110

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1235
$ws4.Range("F6").Value  = 23
$ws4.Range("F7").Value  = 12538
$ws4.Range("F9").Value  = 28
$ws4.Range("F10").Value = 28
$ws4.Range("F12").Value = 12417
$ws4.Range("F13").Value = 240
$ws4.Range("F14").Value = 4896
$ws4.Range("F15").Value = 4823
$ws4.Range("F16").Value = 153
$ws4.Range("F18").Value = 426
$ws4.Range("F19").Value = 110
